$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append Week 15 simulation numbers to the running play-by-play
# sequences (shared strings holding long space-separated number lists).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value2 + " -1 4 0 1 6 5 -6 2 10 0 2 2 -1 5 2 9 0 -3 0 8 3 5 -4 7 4 3"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 7 17 13 30 5 10 7 0 15 8 11 15 15 12 6 9 18 7 11 -2 19 10 7 11 9 11 12 4 17 6 18 7 2 6 8"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 9 13 -1 3 -3 5 3 32 7 -3 8 5 9 3 13 -1 3 13 3 0 3 6 3 47 4 3 0 4 2"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " 1 29 5 8 24 25 16 55 7 13 22 14 1 29 6 1 4"

# ---------------------------------------------------------------------------
# OFF sheet: updated season-aggregate totals after Week 15.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 271
$ws.Range("D2").Value = 12
$ws.Range("F2").Value = 115
$ws.Range("G2").Value = 78
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 54
$ws.Range("L2").Value = 524
$ws.Range("M2").Value = 360
$ws.Range("O2").Value = 35
$ws.Range("Q2").Value = 899

$ws.Range("C3").Value = 304
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 56
$ws.Range("F3").Value = 158
$ws.Range("G3").Value = 55
$ws.Range("H3").Value = 54
$ws.Range("I3").Value = 98
$ws.Range("J3").Value = 107
$ws.Range("N3").Value = 45

# ---------------------------------------------------------------------------
# DEF sheet: updated season-aggregate totals after Week 15.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value = 404
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 107
$ws.Range("G2").Value = 116
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 57
$ws.Range("L2").Value = 489
$ws.Range("M2").Value = 329
$ws.Range("Q2").Value = 970

$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 287
$ws.Range("E3").Value = 52
$ws.Range("F3").Value = 183
$ws.Range("G3").Value = 55
$ws.Range("I3").Value = 94
$ws.Range("J3").Value = 89

# ---------------------------------------------------------------------------
# ST sheet: updated season-aggregate totals, plus appended Week 15 per-kick
# distance sequences (shared strings).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 126
$ws.Range("D2").Value = 125
$ws.Range("F2").Value = 174
$ws.Range("G2").Value = 156
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 70
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1

$ws.Range("B3").Value = 75

$ws.Range("B4").Value = $ws.Range("B4").Value2 + " 55 57 61"
$ws.Range("B5").Value = $ws.Range("B5").Value2 + " 15 18 25"
$ws.Range("B6").Value = $ws.Range("B6").Value2 + " 13 24 11"
$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 54 47 58 35 35"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " 0 0 0 9 22"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 0 2 0 0"

# ---------------------------------------------------------------------------
# TURNS sheet: updated season-aggregate totals.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("E2").Value = 10
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 9

# ---------------------------------------------------------------------------
# PEN sheet: updated season-aggregate totals.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B2").Value = 30
$ws.Range("B3").Value = 37
$ws.Range("D4").Value = 12
